$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the cached "datetimeFigureOut" date field text (20.10.2025 ->
#    24.10.2025) on every Date Placeholder shape found on the slide master
#    and on each of its custom layouts.
# ---------------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($si = 1; $si -le $shapes.Count; $si++) {
        $sh = $shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq "20.10.2025") {
                $sh.TextFrame.TextRange.Text = "24.10.2025"
            }
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2. Remove the first four slides (Architecture pipeline, Structure of NN,
#    Pytorch, Problems to solve), keeping the remaining three slides
#    (Added drop-down menu., Customizable threshold, module architecture)
#    which become slides 1-3.
# ---------------------------------------------------------------------------
for ($i = 4; $i -ge 1; $i--) {
    $p.Slides.Item($i).Delete()
}
